$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.327.95'
$ws.Range('E2').Value = '  +0.70%  '
$ws.Range('D3').Value = '2.281.91'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '499.78'
$ws.Range('E5').Value = '  +1.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.03'
$ws.Range('E6').Value = '  +0.89%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0953'
$ws.Range('E9').Value = '  +1.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.151'
$ws.Range('E10').Value = '  +0.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.333'
$ws.Range('E11').Value = '  +3.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.71'
$ws.Range('E12').Value = '  +1.06%  '
$ws.Range('D13').Value = '2.689.83'
$ws.Range('E13').Value = '  +0.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.88'
$ws.Range('E14').Value = '  +6.19%  '
$ws.Range('D15').Value = '54.289.39'
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('D17').Value = '2.278.65'
$ws.Range('E17').Value = '  +0.51%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.25'
$ws.Range('E18').Value = '  +2.98%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.12'
$ws.Range('E19').Value = '  +2.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '304.53'
$ws.Range('E20').Value = '  +1.79%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.40'
$ws.Range('E21').Value = '  +2.07%  '
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('E23').Value = '  -2.33%  '
$ws.Range('E24').Value = '  -0.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.151'
$ws.Range('E25').Value = '  +2.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.32'
$ws.Range('E26').Value = '  +2.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '174.07'
$ws.Range('E27').Value = '  +6.83%  '
$ws.Range('E28').Value = '  +0.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.98'
$ws.Range('E29').Value = '  +2.56%  '
$ws.Range('D30').Value = '0.0₃0689'
$ws.Range('E30').Value = '  +0.64%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.08'
$ws.Range('E31').Value = '  +1.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '17.80'
$ws.Range('E33').Value = '  +1.70%  '
$ws.Range('E34').Value = '  +0.16%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.933'
$ws.Range('E35').Value = '  +8.39%  '
$ws.Range('E36').Value = '  +1.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.76'
$ws.Range('E37').Value = '  +3.05%  '
$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.374'
$ws.Range('E38').Value = '  -0.34%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.41'
$ws.Range('E39').Value = '  +1.03%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.39'
$ws.Range('E40').Value = '  +1.51%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '124.92'
$ws.Range('E41').Value = '  -1.51%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.78'
$ws.Range('E42').Value = '  -2.16%  '
$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0495'
$ws.Range('E43').Value = '  +3.05%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0896'
$ws.Range('E44').Value = '  +0.58%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.548'
$ws.Range('E45').Value = '  +0.20%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '240.51'
$ws.Range('E46').Value = '  +0.63%  '
$ws.Range('B47').Value = 'Polygon'
$ws.Range('C47').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.372'
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0206'
$ws.Range('E48').Value = '  +1.56%  '
$ws.Range('B49').Value = 'WhiteBITCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.77'
$ws.Range('E49').Value = '  +1.12%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.35'
$ws.Range('E50').Value = '  +0.35%  '
$ws.Range('B51').Value = 'ZEEBU'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.64'
$ws.Range('E51').Value = '  +0.49%  '
